# Case and Fatality Demographics Data Updated
# Updates the underlying counts on each of the six demographic breakdown
# sheets (the week's line-list pull moved from 5.13.21 to 5.20.21), and
# moves the active selection on each sheet as a side effect of refreshing
# the report, ending with "Cases by Age Group" as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Cases by Age Group"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value2  = 278
$ws1.Range("B3").Value2  = 1373
$ws1.Range("B4").Value2  = 3775
$ws1.Range("B5").Value2  = 15743
$ws1.Range("B6").Value2  = 17267
$ws1.Range("B7").Value2  = 15148
$ws1.Range("B8").Value2  = 12785
$ws1.Range("B9").Value2  = 4626
$ws1.Range("B10").Value2 = 3132
$ws1.Range("B11").Value2 = 1900
$ws1.Range("B12").Value2 = 1251
$ws1.Range("B13").Value2 = 1948

# ---------------------------------------------------------------------
# Sheet 2: "Cases by Gender"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value2 = 26953
$ws2.Range("B3").Value2 = 51387
$ws2.Range("B4").Value2 = 899

# ---------------------------------------------------------------------
# Sheet 3: "Cases by RaceEthnicity"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B3").Value2 = 13013
$ws3.Range("B4").Value2 = 28372
$ws3.Range("B5").Value2 = 565
$ws3.Range("B6").Value2 = 27678
$ws3.Range("B7").Value2 = 8655

# ---------------------------------------------------------------------
# Sheet 4: "Fatalities by Age Group"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B5").Value2  = 254
$ws4.Range("B6").Value2  = 850
$ws4.Range("B7").Value2  = 2481
$ws4.Range("B8").Value2  = 5646
$ws4.Range("B9").Value2  = 4723
$ws4.Range("B10").Value2 = 6082
$ws4.Range("B11").Value2 = 6710
$ws4.Range("B12").Value2 = 6615
$ws4.Range("B13").Value2 = 16681

# ---------------------------------------------------------------------
# Sheet 5: "Fatalities by Gender"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value2 = 21015
$ws5.Range("B3").Value2 = 29078

# ---------------------------------------------------------------------
# Sheet 6: "Fatalities by Race-Ethnicity"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B2").Value2 = 1063
$ws6.Range("B3").Value2 = 5062
$ws6.Range("B4").Value2 = 23282
$ws6.Range("B5").Value2 = 273
$ws6.Range("B6").Value2 = 20391

# ---------------------------------------------------------------------
# Refresh the on-screen selection on every sheet (as happens naturally
# while reviewing/updating each tab), finishing on sheet 1 so it is the
# one left active/selected when the workbook is saved.
# ---------------------------------------------------------------------
$ws6.Activate()
$ws6.Range("G14").Select() | Out-Null

$ws5.Activate()
$ws5.Range("D14").Select() | Out-Null

$ws4.Activate()
$ws4.Range("F11").Select() | Out-Null

$ws3.Activate()
$ws3.Range("C14").Select() | Out-Null

$ws2.Activate()
$ws2.Range("E21").Select() | Out-Null

$ws1.Activate()
$ws1.Range("E7").Select() | Out-Null
